# Apply updated odds/value figures to rows 2-5 of Sheet1, matching the
# target commit's data refresh ("Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.8    # was 2.75
$ws.Range("I2").Value = 2.38   # was 2.4
$ws.Range("J2").Value = 3.4    # was 3.25
$ws.Range("K2").Value = 2.2    # was 2.25
$ws.Range("N2").Value = 12     # was 13
$ws.Range("Q2").Value = 1.83   # was 1.73
$ws.Range("R2").Value = 2.03   # was 2.1
$ws.Range("S2").Value = 1.36   # was 1.33
$ws.Range("T2").Value = 3      # was 3.25
$ws.Range("U2").Value = 1.67   # was 1.62
$ws.Range("V2").Value = 2.1    # was 2.2
$ws.Range("AB2").Value = 29    # was 26
$ws.Range("AC2").Value = 12    # was 13
$ws.Range("AD2").Value = 6.5   # was 7
$ws.Range("AI2").Value = 12    # was 13
$ws.Range("AL2").Value = 19    # was 17
$ws.Range("AM2").Value = 26    # was 23
$ws.Range("AP2").Value = 23    # was 21
$ws.Range("AS2").Value = 151   # was 126
$ws.Range("AT2").Value = 3     # was 3.25
$ws.Range("AV2").Value = 51    # was 41
$ws.Range("BC2").Value = 151   # was 126

# Row 3
$ws.Range("G3").Value = 2.6    # was 2.7
$ws.Range("H3").Value = 3.25   # was 3.2
$ws.Range("I3").Value = 2.5    # was 2.45
$ws.Range("J3").Value = 3.4    # was 3.5
$ws.Range("O3").Value = 1.36   # was 1.33
$ws.Range("P3").Value = 3      # was 3.25
$ws.Range("Q3").Value = 2.15   # was 2.08
$ws.Range("R3").Value = 1.67   # was 1.73
$ws.Range("W3").Value = 8      # was 8.5
$ws.Range("Y3").Value = 10     # was 11
$ws.Range("Z3").Value = 26     # was 29
$ws.Range("AC3").Value = 9.5   # was 9
$ws.Range("AD3").Value = 6.5   # was 6
$ws.Range("AG3").Value = 301   # was 251
$ws.Range("AK3").Value = 26    # was 23
$ws.Range("AO3").Value = 15    # was 17

# Row 4
$ws.Range("O4").Value = 1.29   # was 1.25
$ws.Range("P4").Value = 3.5    # was 3.75
$ws.Range("Q4").Value = 1.93   # was 1.9
$ws.Range("R4").Value = 1.93   # was 1.95

# Row 5
$ws.Range("G5").Value = 2.63   # was 2.6
$ws.Range("I5").Value = 2.55   # was 2.6
$ws.Range("L5").Value = 3.25   # was 3.4
$ws.Range("O5").Value = 1.36   # was 1.33
$ws.Range("P5").Value = 3      # was 3.25
$ws.Range("Q5").Value = 2.1    # was 2.05
$ws.Range("R5").Value = 1.7    # was 1.75
$ws.Range("AC5").Value = 9.5   # was 9
$ws.Range("AI5").Value = 12    # was 13
$ws.Range("AM5").Value = 29    # was 34
$ws.Range("AN5").Value = 4.75  # was 4.5
$ws.Range("AZ5").Value = 23    # was 26
$ws.Range("BC5").Value = 151   # was 201
